$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price for "motor drivers" row (B11): 19 -> 38
$ws.Range("B11").Value = 38

# Recalculate dependent formulas (F2 = SUM(B2:B40), F3 = F2/6)
$excel.Calculate()

# Update the active cell selection to C12
$ws.Range("C12").Select()
